$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 98, shifting existing rows 98-118 down to 99-119
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new data record
$ws.Range("A98").Value = 6
$ws.Range("B98").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C98").Value = "Metropolitana"
$ws.Range("D98").Value = 44711
$ws.Range("D98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E98").Value = 13
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100104
$ws.Range("H98").Value = "Frutos de pepita"
$ws.Range("I98").Value = 100104003
$ws.Range("J98").Value = "Membrillo"
$ws.Range("K98").Value = "Champion"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 20
$ws.Range("N98").Value = 220000
$ws.Range("O98").Value = 220000
$ws.Range("P98").Value = 220000
$ws.Range("Q98").Value = "$/bins (400 kilos)"
$ws.Range("R98").Value = "Provincia de Cachapoal"
$ws.Range("S98").Value = 550
$ws.Range("T98").Value = 400
